# game_06_abyssal_nightfall/weapons.xlsx - localize weapon copy + shuffle fx/sfx/notes columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

$ws.Range('S4').Value = 'string'
$ws.Range('T4').Value = 'string'
$ws.Range('V4').Value = 'float'
$ws.Range('W4').Value = 'float'
$ws.Range('S5').Value = 'notes'
$ws.Range('T5').Value = 'fireSfx'
$ws.Range('U5').Value = 'impactSfx'
$ws.Range('V5').Value = 'projectileScale'
$ws.Range('W5').Value = 'impactScale'
$ws.Range('D6').Value = '霓虹符文左轮'
$ws.Range('E6').Value = '霓虹侧臂'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '62'
$ws.Range('G6').ClearFormats()
$ws.Range('P6').Value = 'fx/projectiles/revolver.png'
$ws.Range('Q6').Value = 'fx/impact/sparks.png'
$ws.Range('R6').Value = 'fx/muzzle/flame.png'
$ws.Range('S6').Value = '轻量侧臂，暴击可震慑敌人 0.3 秒。'
$ws.Range('T6').Value = 'ui/assets/sfx/weapons/runic_revolver_fire.wav'
$ws.Range('U6').Value = 'ui/assets/sfx/weapons/runic_revolver_hit.wav'
$ws.Range('V6').NumberFormat = '@'
$ws.Range('V6').Value = '0.68'
$ws.Range('V6').ClearFormats()
$ws.Range('W6').NumberFormat = '@'
$ws.Range('W6').Value = '0.9'
$ws.Range('W6').ClearFormats()
$ws.Range('D7').Value = '共鸣光谱射线'
$ws.Range('E7').Value = '以太光束'
$ws.Range('P7').Value = 'fx/projectiles/beam_ray.png'
$ws.Range('Q7').Value = 'fx/impact/void_burst.png'
$ws.Range('R7').Value = 'fx/muzzle/chorus.png'
$ws.Range('S7').Value = '持续瞄准时，伤害每秒提升 +4。'
$ws.Range('T7').Value = 'ui/assets/sfx/weapons/chorus_ray_fire.wav'
$ws.Range('U7').Value = 'ui/assets/sfx/weapons/chorus_ray_hit.wav'
$ws.Range('V7').NumberFormat = '@'
$ws.Range('V7').Value = '0.85'
$ws.Range('V7').ClearFormats()
$ws.Range('W7').NumberFormat = '@'
$ws.Range('W7').Value = '1.05'
$ws.Range('W7').ClearFormats()
$ws.Range('D8').Value = '潮裂破晓炮'
$ws.Range('E8').Value = '潮汐重炮'
$ws.Range('P8').Value = 'fx/projectiles/tide_shell.png'
$ws.Range('Q8').Value = 'fx/impact/frost_shatter.png'
$ws.Range('R8').Value = 'fx/muzzle/water.png'
$ws.Range('S8').Value = '0.35 秒后裂解成三枚霜晶碎片。'
$ws.Range('T8').Value = 'ui/assets/sfx/weapons/tidebreaker_launcher_fire.wav'
$ws.Range('U8').Value = 'ui/assets/sfx/weapons/tidebreaker_launcher_hit.wav'
$ws.Range('V8').NumberFormat = '@'
$ws.Range('V8').Value = '0.75'
$ws.Range('V8').ClearFormats()
$ws.Range('W8').NumberFormat = '@'
$ws.Range('W8').Value = '1.15'
$ws.Range('W8').ClearFormats()
$ws.Range('D9').Value = '脉冲疾能卡宾枪'
$ws.Range('E9').Value = '脉冲步枪'
$ws.Range('P9').Value = 'fx/projectiles/pulse.png'
$ws.Range('Q9').Value = 'fx/impact/pulse_flash.png'
$ws.Range('R9').Value = 'fx/muzzle/pulse_muzzle.png'
$ws.Range('S9').Value = '跟随节奏射击可叠加动量，加速装填。'
$ws.Range('T9').Value = 'ui/assets/sfx/weapons/pulse_carbine_fire.wav'
$ws.Range('U9').Value = 'ui/assets/sfx/weapons/pulse_carbine_hit.wav'
$ws.Range('V9').NumberFormat = '@'
$ws.Range('V9').Value = '0.7'
$ws.Range('V9').ClearFormats()
$ws.Range('W9').NumberFormat = '@'
$ws.Range('W9').Value = '1.0'
$ws.Range('W9').ClearFormats()
$ws.Range('D10').Value = '幽幕散裂霰炮'
$ws.Range('E10').Value = '暗影霰弹'
$ws.Range('P10').Value = 'fx/projectiles/umbral_pellet.png'
$ws.Range('Q10').Value = 'fx/impact/void_scar.png'
$ws.Range('R10').Value = 'fx/muzzle/umbral.png'
$ws.Range('S10').Value = '近距离命中同一目标时，每颗弹丸额外 +6 伤害。'
$ws.Range('T10').Value = 'ui/assets/sfx/weapons/umbral_scattergun_fire.wav'
$ws.Range('U10').Value = 'ui/assets/sfx/weapons/umbral_scattergun_hit.wav'
$ws.Range('V10').NumberFormat = '@'
$ws.Range('V10').Value = '0.78'
$ws.Range('V10').ClearFormats()
$ws.Range('W10').NumberFormat = '@'
$ws.Range('W10').Value = '1.1'
$ws.Range('W10').ClearFormats()
$ws.Range('D11').Value = '蚀光穿梭矛阵'
$ws.Range('E11').Value = '熔蚀矛阵'
$ws.Range('P11').Value = 'fx/projectiles/eclipse_javelin.png'
$ws.Range('Q11').Value = 'fx/impact/eclipse_burst.png'
$ws.Range('R11').Value = 'fx/muzzle/eclipse.png'
$ws.Range('S11').Value = '自导矛头炸裂出炽焰柱，擅长清理聚群敌人。'
$ws.Range('T11').Value = 'ui/assets/sfx/weapons/eclipse_javelin_fire.wav'
$ws.Range('U11').Value = 'ui/assets/sfx/weapons/eclipse_javelin_hit.wav'
$ws.Range('V11').NumberFormat = '@'
$ws.Range('V11').Value = '0.82'
$ws.Range('V11').ClearFormats()
$ws.Range('W11').NumberFormat = '@'
$ws.Range('W11').Value = '1.25'
$ws.Range('W11').ClearFormats()
